# prescriptions.xlsx: BEGIN_TIME / END_TIME (E2:F2) now come from a
# LocalDateTime-backed column instead of a sql.Time string literal, so the
# sheet should hold a real Excel datetime serial with a custom
# yyyy/mm/dd hh:mm:ss.ss display format rather than the old shared-string
# timestamp text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the new custom number format first so the upcoming date values pick
# it up directly (matches the workbook's existing custom-numFmt reuse
# pattern instead of bouncing through a default date/time format).
$ws.Range("E2:F2").NumberFormat = "yyyy/mm/dd\ hh:mm:ss\.ss"

$dt = Get-Date -Year 2022 -Month 4 -Day 22 -Hour 10 -Minute 34 -Second 0
$ws.Range("E2").Value = $dt
$ws.Range("F2").Value = $dt

# The author's last selection in the sheet ended up on E2 (BEGIN_TIME).
$ws.Range("E2").Select()
